$d = $word.ActiveDocument

$replacements = @(
    @("569÷4=", "508÷8="),
    @("965÷7=", "469÷5="),
    @("449÷7=", "446÷8="),
    @("434÷4=", "259÷7="),
    @("268÷8=", "505÷4="),
    @("285÷4=", "652÷7="),
    @("999÷6=", "877÷3="),
    @("435÷9=", "926÷6="),
    @("981÷9=", "785÷7="),
    @("907÷2=", "644÷6="),
    @("520÷3=", "476÷3="),
    @("105÷6=", "766÷5="),
    @("990÷8=", "288÷4="),
    @("172÷9=", "843÷4="),
    @("955÷9=", "498÷2="),
    @("791÷5=", "220÷4="),
    @("551÷6=", "323÷2="),
    @("853÷3=", "821÷8="),
    @("702÷5=", "182÷4="),
    @("552÷3=", "698÷4="),
    @("259÷3=", "505÷7="),
    @("697÷2=", "948÷4="),
    @("195÷7=", "925÷9="),
    @("861÷9=", "632÷4="),
    @("128÷9=", "918÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
